$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 67.666664
$ws.Range("I11").Value = 67.666664
$ws.Range("K11").Value = 67.666664
$ws.Range("M11").Value = 72.333336
$ws.Range("H17").Value = 3879.1482
$ws.Range("J17").Value = 3879.1482
$ws.Range("L17").Value = 11637.4446
$ws.Range("N17").Value = -11973.4446
$ws.Range("H40").Value = 2197.8667
$ws.Range("I40").Value = 1995.4286
$ws.Range("J40").Value = 2375
$ws.Range("K40").Value = 1995.4286
$ws.Range("L40").Value = 2375
$ws.Range("M40").Value = -1820.4286
$ws.Range("N40").Value = -2725
$ws.Range("H64").Value = 3198
$ws.Range("I64").Value = 3190
$ws.Range("K64").Value = 3190
$ws.Range("M64").Value = -2942
$ws.Range("H67").Value = 3198
$ws.Range("I67").Value = 3190
$ws.Range("K67").Value = 3190
$ws.Range("M67").Value = -2332
$ws.Range("H74").Value = 7990
$ws.Range("I74").Value = 7990
$ws.Range("K74").Value = 7990
$ws.Range("M74").Value = -7054
$ws.Range("H77").Value = 7990
$ws.Range("I77").Value = 7990
$ws.Range("K77").Value = 39950
$ws.Range("M77").Value = -35270
$ws.Range("H86").Value = 4916.3335
$ws.Range("I86").Value = 4750
$ws.Range("J86").Value = 4999.5
$ws.Range("K86").Value = 4750
$ws.Range("L86").Value = 4999.5
$ws.Range("M86").Value = -3627
$ws.Range("N86").Value = -7245.5
$ws.Range("H89").Value = 4916.3335
$ws.Range("I89").Value = 4750
$ws.Range("J89").Value = 4999.5
$ws.Range("K89").Value = 23750
$ws.Range("L89").Value = 24997.5
$ws.Range("M89").Value = -18134
$ws.Range("N89").Value = -36229.5
$ws.Range("H101").Value = 1352.25
$ws.Range("I101").Value = 1352.25
$ws.Range("K101").Value = 4056.75
$ws.Range("M101").Value = -2434.75
$ws.Range("H106").Value = 40384.7
$ws.Range("I106").Value = 40384.7
$ws.Range("K106").Value = 40384.7
$ws.Range("M106").Value = -39753.7
$ws.Range("H112").Value = 2187.2354
$ws.Range("J112").Value = 2187.2354
$ws.Range("L112").Value = 6561.706200000001
$ws.Range("N112").Value = -8777.706200000001
$ws.Range("H113").Value = 2349
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()  # was -9508
$ws.Range("H129").Value = 3767.182
$ws.Range("I129").Value = 2364.8333
$ws.Range("K129").Value = 7094.499899999999
$ws.Range("M129").Value = -2094.499899999999
$ws.Range("H132").Value = 2850.4285
$ws.Range("I132").Value = 2492.1667
$ws.Range("K132").Value = 7476.500100000001
$ws.Range("M132").Value = -4946.500100000001
$ws.Range("H138").Value = 2220.56
$ws.Range("I138").Value = 1329.3334
$ws.Range("J138").Value = 2502
$ws.Range("K138").Value = 3988.0002
$ws.Range("L138").Value = 7506
$ws.Range("M138").Value = 1151.9998
$ws.Range("N138").Value = -17786

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1764.3334
$ws.Range("I45").Value = 1764.3636
$ws.Range("K45").Value = 1764.3636
$ws.Range("M45").Value = -1387.3636
$ws.Range("H61").Value = 1798.2916
$ws.Range("I61").Value = 1622
$ws.Range("K61").Value = 1622
$ws.Range("M61").Value = -1410
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H131").Value = 100000
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H132").Value = 2061.75
$ws.Range("I132").Value = 2046
$ws.Range("K132").Value = 6138
$ws.Range("M132").Value = -3608
$ws.Range("H136").Value = 1798.2916
$ws.Range("I136").Value = 1622
$ws.Range("K136").Value = 4866
$ws.Range("M136").Value = -2316

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6807
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()  # was -872
$ws.Range("H89").Value = 6807
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()  # was -4359
$ws.Range("H94").Value = 969.4375
$ws.Range("I94").Value = 1005.3333
$ws.Range("K94").Value = 1005.3333
$ws.Range("M94").Value = -554.3333
$ws.Range("H99").Value = 1439.5428
$ws.Range("I99").Value = 1302.7037
$ws.Range("K99").Value = 1302.7037
$ws.Range("M99").Value = 195.2963
$ws.Range("H105").Value = 2943.9375
$ws.Range("I105").Value = 2848.4285
$ws.Range("K105").Value = 2848.4285
$ws.Range("M105").Value = -1101.4285
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H134").Value = 2466.9524
$ws.Range("I134").Value = 2179.9285
$ws.Range("K134").Value = 6539.7855
$ws.Range("M134").Value = -4004.7855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3818.024
$ws.Range("I31").Value = 3237.0557
$ws.Range("K31").Value = 3237.0557
$ws.Range("M31").Value = -2942.0557
$ws.Range("H34").Value = 3818.024
$ws.Range("I34").Value = 3237.0557
$ws.Range("K34").Value = 3237.0557
$ws.Range("M34").Value = -3035.0557
$ws.Range("H58").Value = 2380.739
$ws.Range("I58").Value = 1162.1666
$ws.Range("K58").Value = 1162.1666
$ws.Range("M58").Value = -959.1666
$ws.Range("H62").Value = 84495.8
$ws.Range("I62").Value = 5163.3335
$ws.Range("J62").Value = 203494.5
$ws.Range("K62").Value = 5163.3335
$ws.Range("L62").Value = 203494.5
$ws.Range("M62").Value = -4539.3335
$ws.Range("N62").Value = -204742.5
$ws.Range("H65").Value = 84495.8
$ws.Range("I65").Value = 5163.3335
$ws.Range("J65").Value = 203494.5
$ws.Range("K65").Value = 25816.6675
$ws.Range("L65").Value = 1017472.5
$ws.Range("M65").Value = -22696.6675
$ws.Range("N65").Value = -1023712.5
$ws.Range("H132").Value = 1737.6271
$ws.Range("J132").Value = 4515.4
$ws.Range("L132").Value = 13546.2
$ws.Range("N132").Value = -18606.2
$ws.Range("H134").Value = 2469.4883
$ws.Range("I134").Value = 2055.516
$ws.Range("J134").Value = 3538.9167
$ws.Range("K134").Value = 6166.548000000001
$ws.Range("L134").Value = 10616.7501
$ws.Range("M134").Value = -3631.548000000001
$ws.Range("N134").Value = -15686.7501
$ws.Range("H136").Value = 2380.739
$ws.Range("I136").Value = 1162.1666
$ws.Range("K136").Value = 3486.4998
$ws.Range("M136").Value = -936.4998000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2133
$ws.Range("I80").Value = 2133
$ws.Range("K80").Value = 6399
$ws.Range("M80").Value = -5463
$ws.Range("H83").Value = 2133
$ws.Range("I83").Value = 2133
$ws.Range("K83").Value = 19197
$ws.Range("M83").Value = -14517
$ws.Range("H136").Value = 14062
$ws.Range("I136").Value = 9608.5
$ws.Range("K136").Value = 28825.5
$ws.Range("M136").Value = -23725.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2186.7222
$ws.Range("I132").Value = 1765
$ws.Range("J132").Value = 2608.4443
$ws.Range("K132").Value = 5295
$ws.Range("L132").Value = 7825.3329
$ws.Range("M132").Value = -2765
$ws.Range("N132").Value = -12885.3329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 497
$ws.Range("I55").Value = 497
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 497
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -324
$ws.Range("N55").ClearContents()  # was -841
$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -751
$ws.Range("N68").ClearContents()  # was -2798
$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -3756
$ws.Range("N71").ClearContents()  # was -13988
$ws.Range("H82").Value = 1963.4615
$ws.Range("I82").Value = 2032.6
$ws.Range("J82").Value = 1733
$ws.Range("K82").Value = 2032.6
$ws.Range("L82").Value = 1733
$ws.Range("M82").Value = -1671.6
$ws.Range("N82").Value = -2455
$ws.Range("H85").Value = 1963.4615
$ws.Range("I85").Value = 2032.6
$ws.Range("J85").Value = 1733
$ws.Range("K85").Value = 2032.6
$ws.Range("L85").Value = 1733
$ws.Range("M85").Value = -784.5999999999999
$ws.Range("N85").Value = -4229
$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752
$ws.Range("H132").Value = 4555.7393
$ws.Range("I132").Value = 3483.3076
$ws.Range("K132").Value = 10449.9228
$ws.Range("M132").Value = -7919.9228
$ws.Range("H136").Value = 4629.5
$ws.Range("I136").Value = 3996
$ws.Range("K136").Value = 11988
$ws.Range("M136").Value = -9438

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 124500
$ws.Range("J117").Value = 124500
$ws.Range("L117").Value = 124500
$ws.Range("N117").Value = -133678
